# Commit: "with IFLT pngs added, folder structure updated"
#
# This script applies the portion of the change that is reachable through
# the Excel object model exposed by this COM-interop host.  The workbook's
# data worksheet ("Sheet" -> "Data") gains two new columns (K: AVG_IND(4),
# L: BEST_IND(3)) that back the updated "Best & Average Individuals" chart,
# and the worksheet itself is renamed to match the new folder/sheet naming
# convention described in the commit message.
#
# NOTE: the two chartsheets ("Math, with IFLT" / "Chart1") and their chart
# parts are not modelled as addressable Worksheets/Charts objects in this
# host (Workbook.Worksheets.Count / Workbook.Charts.Count == 1 and 0
# respectively for this file), so their titles/series/axes cannot be
# touched from here; only the underlying worksheet data model is editable.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- rename the data sheet: "Sheet" -> "Data" -----------------------------
$ws.Name = "Data"

# --- new header cells: K1 = AVG_IND(4), L1 = BEST_IND(3) ------------------
$ws.Cells.Item(1, 11).Value = "AVG_IND(4)"
$ws.Cells.Item(1, 12).Value = "BEST_IND(3)"

# --- new data columns K2:L79 ----------------------------------------------
$klData = @(
    @(2, 0.95104327359190199, 0.91568111042655498),
    @(3, 0.95103342437482097, 0.915762463107887),
    @(4, 0.95101025252284899, 0.91584360648490304),
    @(5, 0.95097420969856405, 0.91592454169376603),
    @(6, 0.95092607678337904, 0.91600526986150099),
    @(7, 0.95086697210078697, 0.91608579210607399),
    @(8, 0.95079835509871102, 0.91616610953647604),
    @(9, 0.95072202098536496, 0.91624622325279803),
    @(10, 0.95064008163269498, 0.91632613434631105),
    @(11, 0.95055492856791901, 0.91640584389954305),
    @(12, 0.950469175109797, 0.91648535298635803),
    @(13, 0.95038557662253098, 0.91656466267202896),
    @(14, 0.95030693029874203, 0.91664377401331598),
    @(15, 0.95023595857659804, 0.91672268805854196),
    @(16, 0.95017518289101299, 0.91680140584766701),
    @(17, 0.95012679656009602, 0.91687992841236199),
    @(18, 0.95129136085021804, 0.916958256776083),
    @(19, 0.95127000910613302, 0.91703639195414399),
    @(20, 0.95126445734965503, 0.91711433495379102),
    @(21, 0.95008353239292498, 0.917192086774272),
    @(22, 0.95011156540348696, 0.91726964840691105),
    @(23, 0.95015343405362995, 0.91734702083517805),
    @(24, 0.95020728881981098, 0.91742420503475997),
    @(25, 0.95027084977492704, 0.91750120197362905),
    @(26, 0.95034152587626597, 0.91757801261211702),
    @(27, 0.95041654172258805, 0.91765463790297996),
    @(28, 0.95049306194343797, 0.91773107879146898),
    @(29, 0.95056830481545795, 0.91780733621539601),
    @(30, 0.95063963893744596, 0.91788341110520799),
    @(31, 0.95070465945555205, 0.91795930438404505),
    @(32, 0.95076124301792697, 0.91803501696781498),
    @(33, 0.95080758300944401, 0.91811054976525297),
    @(34, 0.95084220841851996, 0.91818590367799402),
    @(35, 0.950863990775899, 0.91826107960063097),
    @(36, 0.95087214394222597, 0.91833607842078702),
    @(37, 0.95086622115833697, 0.918410901019173),
    @(38, 0.95084611282411202, 0.91848554826965401),
    @(39, 0.95081204709328704, 0.918560021039314),
    @(40, 0.95076459373947197, 0.91863432018851499),
    @(41, 0.95070467005002401, 0.91870844657096296),
    @(42, 0.95063354593068095, 0.91878240103376696),
    @(43, 0.95055284414514396, 0.91885618441750105),
    @(44, 0.95046453085023397, 0.91892979755626603),
    @(45, 0.95037089147324505, 0.91900324127774702),
    @(46, 0.95027448761986399, 0.91907651640327803),
    @(47, 0.95017809212803195, 0.919149623747894),
    @(48, 0.95008460152306096, 0.91922256412039605),
    @(49, 0.94999692779142497, 0.91929533832340704),
    @(50, 0.94991787426918695, 0.91936794715342796),
    @(51, 0.94985000314347501, 0.91944039140089695),
    @(52, 0.94979550416566805, 0.91951267185024699),
    @(53, 0.95083386223808497, 0.91958478927995901),
    @(54, 0.95080641835940105, 0.91965674446261902),
    @(55, 0.95079551377898797, 0.91972853816497402),
    @(56, 0.94973600432511596, 0.91980017114798596),
    @(57, 0.94976131163461397, 0.91987164416688705),
    @(58, 0.94980062249436303, 0.91994295797123105),
    @(59, 0.94985189766161904, 0.92001411330494998),
    @(60, 0.94991268183592303, 0.92008511090640199),
    @(61, 0.94998023242873497, 0.92015595150842999),
    @(62, 0.95005165334888197, 0.92022663583840703),
    @(63, 0.95012402360400505, 0.92029716461829303),
    @(64, 0.95019451230234997, 0.92036753856468001),
    @(65, 0.95026047419156101, 0.920437758388849),
    @(66, 0.95031952276713205, 0.92050782479681503),
    @(67, 0.95036958080422096, 0.92057773848937596),
    @(68, 0.95040891056292898, 0.92064750016216601),
    @(69, 0.95043612764733099, 0.92071711050570004),
    @(70, 0.95045020344216502, 0.92078657020542398),
    @(71, 0.95045046119952503, 0.92085587994175999),
    @(72, 0.95043657028123996, 0.92092504039015699),
    @(73, 0.95040854192129998, 0.92099405222113295),
    @(74, 0.95036672833242797, 0.92106291610032498),
    @(75, 0.95031182523426305, 0.92113163268853204),
    @(76, 0.95024487612644504, 0.92120020264176306),
    @(77, 0.95016727506719001, 0.92126862661127695),
    @(78, 0.95008076354061399, 0.92133690524363498),
    @(79, 0.94998741638216899, 0.92140503918073502)
)

foreach ($row in $klData) {
    $r = $row[0]
    $ws.Cells.Item($r, 11).Value = $row[1]
    $ws.Cells.Item($r, 12).Value = $row[2]
}

# --- best-effort view/selection update -------------------------------------
# The target selection spans two disjoint areas (A1:A79 and K1:L79) with
# A1:A79 active-cell context pointing at K1; this host's Selection model
# only keeps the most recently selected single area, so we select the new
# range last to land on activeCell = K1 as closely as this model allows.
$ws.Range("A1:A79").Select()
$ws.Range("K1:L79").Select()
